$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 615-616 (pushing the existing 2026/12/29 block and
# everything below it down by two rows), then populate them with the new
# 2026/01/09 (23:00) and 2026/01/10 (02:00) entries.
$ws.Range("A615:D616").EntireRow.Insert()

# Row 615: 2026/01/09, 金, 23, 201
$ws.Range("A615").Value = "'2026/01/09"
$ws.Range("A615").Style = "Normal"
$ws.Range("B615").Value = "金"
$ws.Range("C615").Value = 23
$ws.Range("D615").Value = 201

# Row 616: 2026/01/10, 土, 2, 201
$ws.Range("A616").Value = "'2026/01/10"
$ws.Range("A616").Style = "Normal"
$ws.Range("B616").Value = "土"
$ws.Range("C616").Value = 2
$ws.Range("D616").Value = 201
